$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '41.900.75'
$ws.Range("E2").Value = '  +6.08%  '
$ws.Range("D3").Value = '2.260.66'
$ws.Range("E3").Value = '  +4.29%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").Value = "'234.47"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.39%  '
$ws.Range("D6").Value = "'0.644"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.01%  '
$ws.Range("D7").Value = "'63.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").Value = "'0.410"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +3.73%  '
$ws.Range("D10").Value = "'59.80"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.92%  '
$ws.Range("D11").Value = "'0.0898"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.35%  '
$ws.Range("E12").Value = '  +1.74%  '
$ws.Range("D13").Value = '2.596.35'
$ws.Range("E13").Value = '  +4.42%  '
$ws.Range("D14").Value = "'16.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '
$ws.Range("D15").Value = "'22.85"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.67%  '
$ws.Range("E16").Value = '  +1.20%  '
$ws.Range("E17").Value = '  +3.45%  '
$ws.Range("D18").Value = '2.260.11'
$ws.Range("E18").Value = '  +4.46%  '
$ws.Range("D19").Value = '41.699.95'
$ws.Range("E19").Value = '  +5.68%  '
$ws.Range("D20").Value = "'74.81"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +4.03%  '
$ws.Range("E21").Value = '  +9.81%  '
$ws.Range("E22").Value = '  -1.15%  '
$ws.Range("D23").Value = "'252.70"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.07%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").Value = "'2.42"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.44%  '
$ws.Range("E26").Value = '  +0.82%  '
$ws.Range("E27").Value = '  +8.25%  '
$ws.Range("D28").Value = "'9.84"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.97%  '
$ws.Range("D29").Value = "'170.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.75%  '
$ws.Range("D30").Value = "'20.51"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.07%  '
$ws.Range("E31").Value = '  +2.23%  '
$ws.Range("E32").Value = '  +7.91%  '
$ws.Range("E33").Value = '  +1.90%  '
$ws.Range("D34").Value = "'5.13"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +8.65%  '
$ws.Range("D35").Value = "'4.79"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.62%  '
$ws.Range("D36").Value = "'0.0641"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.39%  '
$ws.Range("D37").Value = "'6.86"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.90%  '
$ws.Range("D38").Value = "'3.84"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +6.30%  '
$ws.Range("E39").Value = '  +0.16%  '
$ws.Range("D40").Value = "'0.000262"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +56.15%  '
$ws.Range("D41").Value = "'5.15"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +20.12%  '
$ws.Range("D43").Value = "'0.0242"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.87%  '
$ws.Range("E44").Value = '  +11.93%  '
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = "'17.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.26%  '
$ws.Range("B46").Value = 'Aave'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D46").Value = "'102.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.94%  '
$ws.Range("E47").Value = '  +6.58%  '
$ws.Range("E48").Value = '  +3.38%  '
$ws.Range("D49").Value = '1.504.25'
$ws.Range("E49").Value = '  -1.39%  '
$ws.Range("E50").Value = '  +1.41%  '
$ws.Range("D51").Value = "'2.82"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.02%  '
